# "Add static folder and js" — the underlying xlsx edit duplicates the
# existing data block (rows 2-9) twice more (into rows 10-17 and 18-25)
# and refreshes the "uuid" column (G) for every data row to a new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate the data block A2:H9 into A10:H17 ---
# Paste formats first (so date/number styles on D,E,H survive) then paste
# values on top, rather than a single combined paste, which would otherwise
# pick up a brand-new "pasted" number format instead of reusing the
# worksheet's existing style indexes.
$ws.Range("A2:H9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:H9").Copy()
$ws.Range("A10").PasteSpecial(-4163)   # xlPasteValues

# --- Duplicate the data block A2:H9 into A18:H25 ---
$ws.Range("A2:H9").Copy()
$ws.Range("A18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:H9").Copy()
$ws.Range("A18").PasteSpecial(-4163)   # xlPasteValues

# --- Refresh the uuid column for every data row (2-25) ---
$ws.Range("G2:G25").Value = "ddb71f0e-ca75-4b40-9ae2-33afa81c43ba"
